# Update cryptos list values per upstream data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.828.31'
$ws.Range('E2').Value = '  +0.27%  '
$ws.Range('D3').Value = '2.085.03'
$ws.Range('E3').Value = '  +0.37%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '233.59'
$ws.Range('E5').Value = '  +0.41%  '
$ws.Range('E6').Value = '  +0.33%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.92'
$ws.Range('E7').Value = '  +3.35%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  +2.04%  '
$ws.Range('E10').Value = '  +0.71%  '
$ws.Range('E11').Value = '  +2.85%  '
$ws.Range('D12').Value = '2.391.45'
$ws.Range('E12').Value = '  +0.72%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.74'
$ws.Range('E13').Value = '  +2.14%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.23'
$ws.Range('E14').Value = '  +1.49%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.780'
$ws.Range('E15').Value = '  +2.73%  '
$ws.Range('E16').Value = '  +1.87%  '
$ws.Range('D17').Value = '2.078.37'
$ws.Range('E17').Value = '  +0.02%  '
$ws.Range('D18').Value = '37.742.33'
$ws.Range('E18').Value = '  +0.23%  '
$ws.Range('E19').Value = '  -0.01%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.77'
$ws.Range('E20').Value = '  +1.50%  '
$ws.Range('D21').Value = '0.0₃0847'
$ws.Range('E21').Value = '  +3.44%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '228.45'
$ws.Range('E22').Value = '  +0.34%  '
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('E24').Value = '  -0.47%  '
$ws.Range('E25').Value = '  +1.74%  '
$ws.Range('B26').Value = 'Monero'
$ws.Range('C26').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '171.06'
$ws.Range('E26').Value = '  +0.80%  '
$ws.Range('B27').Value = 'Cosmos'
$ws.Range('C27').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.52'
$ws.Range('E27').Value = '  +7.00%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.138'
$ws.Range('E28').Value = '  -0.71%  '
$ws.Range('E29').Value = '  +0.22%  '
$ws.Range('E30').Value = '  +1.20%  '
$ws.Range('E31').Value = '  +2.60%  '
$ws.Range('E32').Value = '  +3.02%  '
$ws.Range('E33').Value = '  +1.80%  '
$ws.Range('E34').Value = '  +2.10%  '
$ws.Range('E35').Value = '  +0.16%  '
$ws.Range('E36').Value = '  +2.02%  '
$ws.Range('E37').Value = '  -0.47%  '
$ws.Range('E38').Value = '  -0.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.43'
$ws.Range('E39').Value = '  +0.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0983'
$ws.Range('E40').Value = '  -0.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '99.27'
$ws.Range('E41').Value = '  +0.62%  '
$ws.Range('E42').Value = '  +3.02%  '
$ws.Range('E43').Value = '  +10.12%  '
$ws.Range('E44').Value = '  -0.69%  '
$ws.Range('D45').Value = '1.452.03'
$ws.Range('E45').Value = '  -0.13%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.15'
$ws.Range('E46').Value = '  -0.39%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '4.21'
$ws.Range('E47').Value = '  -3.29%  '
$ws.Range('E48').Value = '  +1.70%  '
$ws.Range('E49').Value = '  -0.52%  '
$ws.Range('E50').Value = '  -0.42%  '
$ws.Range('D51').Value = '2.277.43'
$ws.Range('E51').Value = '  +0.52%  '
